$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-19) holds the "Förändrad" (changed) date, stored as serial 45183
# (2023-09-14). The automatic update bumps it by one day to 45184 (2023-09-15).
$ws.Range("C2:C19").Value = 45184
